$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModificarClienteCorporativo")

# Update the RUC and "RUC Modificar" values in row 2 (keep them as text,
# matching the existing quotePrefix text formatting of these cells).
$ws.Range("C2").Value = "'82584292"
$ws.Range("D2").Value = "'20552103816"
